$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (A2 unchanged = ECs, D2 changes from M2 to ECs)
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 18.76479233333333
$ws.Range("H2").Value = 56.294377
$ws.Range("I2").Value = 0.09818846546758657
$ws.Range("J2").Value = 0.09818846546758656
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.438062
$ws.Range("N2").Value = 1.314186
$ws.Range("O2").Value = 0.6074000808827777
$ws.Range("P2").Value = 0.6074000808827777
$ws.Range("Q2").Value = 8.220142459124666
$ws.Range("R2").Value = 73.981282132122
$ws.Range("S2").Value = 0.05963968186676791
$ws.Range("T2").Value = 0.0596396818667679

# Row 3: A3 changes from FAPs to ECs, D3 stays M2
$ws.Range("A3").Value = "ECs"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 18.76479233333333
$ws.Range("H3").Value = 56.294377
$ws.Range("I3").Value = 0.09818846546758657
$ws.Range("J3").Value = 0.09818846546758656
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.2831463333333333
$ws.Range("N3").Value = 0.8494390000000001
$ws.Range("O3").Value = 0.3925999191172223
$ws.Range("P3").Value = 0.3925999191172223
$ws.Range("Q3").Value = 5.313182144944777
$ws.Range("R3").Value = 47.818639304503
$ws.Range("S3").Value = 0.03854878360081865
$ws.Range("T3").Value = 0.03854878360081865

# Row 4: A4 changes from M2 to FAPs, D4 changes from M2 to ECs
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 72.11798333333333
$ws.Range("H4").Value = 216.35395
$ws.Range("I4").Value = 0.3773638413007209
$ws.Range("J4").Value = 0.3773638413007209
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.438062
$ws.Range("N4").Value = 1.314186
$ws.Range("O4").Value = 0.6074000808827777
$ws.Range("P4").Value = 0.6074000808827777
$ws.Range("Q4").Value = 31.59214801496666
$ws.Range("R4").Value = 284.3293321347
$ws.Range("S4").Value = 0.2292108277282936
$ws.Range("T4").Value = 0.2292108277282936

# Row 5: A5 changes from sCs to FAPs, D5 stays M2
$ws.Range("A5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 72.11798333333333
$ws.Range("H5").Value = 216.35395
$ws.Range("I5").Value = 0.3773638413007209
$ws.Range("J5").Value = 0.3773638413007209
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.2831463333333333
$ws.Range("N5").Value = 0.8494390000000001
$ws.Range("O5").Value = 0.3925999191172223
$ws.Range("P5").Value = 0.3925999191172223
$ws.Range("Q5").Value = 20.41994254822778
$ws.Range("R5").Value = 183.77948293405
$ws.Range("S5").Value = 0.1481530135724273
$ws.Range("T5").Value = 0.1481530135724273

# Row 6 (new): M2 / Cxcl10 / Ccr3 / ECs
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Cxcl10"
$ws.Range("C6").Value = "Ccr3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 83.31930033333333
$ws.Range("H6").Value = 249.957901
$ws.Range("I6").Value = 0.4359757410707099
$ws.Range("J6").Value = 0.4359757410707098
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.438062
$ws.Range("N6").Value = 1.314186
$ws.Range("O6").Value = 0.6074000808827777
$ws.Range("P6").Value = 0.6074000808827777
$ws.Range("Q6").Value = 36.49901934262066
$ws.Range("R6").Value = 328.491174083586
$ws.Range("S6").Value = 0.2648117003892781
$ws.Range("T6").Value = 0.2648117003892781

# Row 7 (new): M2 / Cxcl10 / Ccr3 / M2
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Cxcl10"
$ws.Range("C7").Value = "Ccr3"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 83.31930033333333
$ws.Range("H7").Value = 249.957901
$ws.Range("I7").Value = 0.4359757410707099
$ws.Range("J7").Value = 0.4359757410707098
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2831463333333333
$ws.Range("N7").Value = 0.8494390000000001
$ws.Range("O7").Value = 0.3925999191172223
$ws.Range("P7").Value = 0.3925999191172223
$ws.Range("Q7").Value = 23.59155438528211
$ws.Range("R7").Value = 212.323989467539
$ws.Range("S7").Value = 0.1711640406814317
$ws.Range("T7").Value = 0.1711640406814317

# Row 8 (new): sCs / Cxcl10 / Ccr3 / ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cxcl10"
$ws.Range("C8").Value = "Ccr3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 16.90787
$ws.Range("H8").Value = 50.72361
$ws.Range("I8").Value = 0.08847195216098278
$ws.Range("J8").Value = 0.08847195216098278
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.438062
$ws.Range("N8").Value = 1.314186
$ws.Range("O8").Value = 0.6074000808827777
$ws.Range("P8").Value = 0.6074000808827777
$ws.Range("Q8").Value = 7.40669534794
$ws.Range("R8").Value = 66.66025813146001
$ws.Range("S8").Value = 0.05373787089843818
$ws.Range("T8").Value = 0.05373787089843818

# Row 9 (new): sCs / Cxcl10 / Ccr3 / M2
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cxcl10"
$ws.Range("C9").Value = "Ccr3"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 16.90787
$ws.Range("H9").Value = 50.72361
$ws.Range("I9").Value = 0.08847195216098278
$ws.Range("J9").Value = 0.08847195216098278
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2831463333333333
$ws.Range("N9").Value = 0.8494390000000001
$ws.Range("O9").Value = 0.3925999191172223
$ws.Range("P9").Value = 0.3925999191172223
$ws.Range("Q9").Value = 4.787401394976667
$ws.Range("R9").Value = 43.08661255479
$ws.Range("S9").Value = 0.0347340812625446
$ws.Range("T9").Value = 0.0347340812625446
